# Insert a new record at row 70 (pushing the existing rows 70-185 down to
# 71-186) and populate it with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(70).Insert()

$ws.Range("A70").Value = 10
$ws.Range("B70").Value = "Vega Modelo de Temuco"
$ws.Range("C70").Value = "La Araucanía"
$ws.Range("D70").Value = 44477
$ws.Range("E70").Value = 9
$ws.Range("F70").Value = 100112017
$ws.Range("G70").Value = "Apio"
$ws.Range("H70").Value = "Americana (o)"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 40
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = 9000
$ws.Range("N70").Value = "$/docena de matas"
$ws.Range("O70").Value = "Provincia del Elquí"
$ws.Range("P70").Value = 1500
$ws.Range("Q70").Value = 6
$ws.Range("R70").Value = "Hortaliza"
